# Update hotel reviews data:
# The "hotel_info" sheet has three columns (English_Reviews_num, Local_Rank,
# Total_Reviews_num) that were left blank for the WoodSpring Suites Houston
# Northwest row. Fill them in with the scraped values. The source data
# stores these (numeric-looking) values as text, same as the neighbouring
# Zip column, so force each cell to text before writing, then drop the
# now-unneeded number format so the cell keeps its default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$ws.Range("G2").Value = "'2"
$ws.Range("G2").ClearFormats()

$ws.Range("H2").Value = "'362"
$ws.Range("H2").ClearFormats()

$ws.Range("I2").Value = "'3"
$ws.Range("I2").ClearFormats()
